$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Oct 02 16:44:59 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 02 16:45:14 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 02 16:45:28 EDT 2023"
